$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The reservation entries below the header are being replaced: only the
# "Libre" (free) rows stay, everything else gets wiped out, and the data
# shrinks down to just two reservation rows.
$ws.Range("A2:H18").EntireRow.Delete()

# Row 2: new reservation entry
$ws.Range("A2").Value = 3
$ws.Range("C2").Value = "Libre"
$ws.Range("D2").Value = "'4"

# Row 3: new reservation entry (Usuario left blank -> "reservar" button slot)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "'"
$ws.Range("C3").Value = "Libre"
$ws.Range("D3").Value = "'4"

$ws.Range("G14").Select()
